$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": reset the previously non-zero sales
# figures back to 0 (automated update wiped out this period's data).
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$zeroedCells = @(
  "E2","N2",
  "D8","M8","O8","R8",
  "D10","M10",
  "D12","L12","M12",
  "D13",
  "D15","M15","O15","P15",
  "E19","H19","I19","L19","M19","P19",
  "M22",
  "O24",
  "D27","L27",
  "M29",
  "M31"
)

foreach ($addr in $zeroedCells) {
  $ws1.Range($addr).Value2 = 0
}

# Row 34 holds the "<n> de 32" counters per product column; since every
# non-zero cell above was reset, every counter collapses to "0 de 32".
$counterCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($c in $counterCols) {
  $ws1.Range($c + "34").Value2 = "0 de 32"
}

# -----------------------------------------------------------------
# Sheet "VENTA MENSUAL": roll the monthly window forward by one month
# (mayo/junio/julio/agosto -> junio/julio/agosto/septiembre), shifting
# each asesor's figures left and bringing in a fresh (empty) month.
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").Value2 = "junio"
$ws2.Range("D1").Value2 = "julio"
$ws2.Range("E1").Value2 = "agosto"
$ws2.Range("F1").Value2 = "septiembre"

for ($r = 2; $r -le 34; $r++) {
  $dVal = $ws2.Range("D$r").Value2
  $eVal = $ws2.Range("E$r").Value2
  $fVal = $ws2.Range("F$r").Value2

  $ws2.Range("C$r").Value2 = $dVal
  $ws2.Range("D$r").Value2 = $eVal
  $ws2.Range("E$r").Value2 = $fVal
  $ws2.Range("F$r").Value2 = 0
}

# Column widths for the (now shifted) month columns D, E and F also
# changed slightly. Excel's ColumnWidth property is offset from the
# raw OOXML "width" attribute by the standard 0.8333... (5px) padding
# for the default Calibri 11 font, so compensate for that offset.
$colWidthOffset = 0.8333333333333334
$ws2.Columns.Item(4).ColumnWidth = 13 - $colWidthOffset
$ws2.Columns.Item(5).ColumnWidth = 14 - $colWidthOffset
$ws2.Columns.Item(6).ColumnWidth = 16 - $colWidthOffset
